# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoffs have now been handed back in sync with en-US:
#   - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#     (Overview sheet + each language sheet)
#   - Each language sheet's rows for the two tracked content files gain a
#     "Latest Target File" (E) and "Latest Handback File" (F) entry (with
#     hyperlinks, mirroring the existing Source File / Handoff File links)
#   - "Latest Handback DateTime" (G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Cornflower blue (FF6495ED) expressed as the BGR-packed OLE_COLOR Excel's
# Font.Color setter expects, so the new hyperlink-look cells match the
# workbook's existing custom "HyperLink" cell style exactly.
$hyperlinkColor = 15570276

function Set-HandbackStatus($ws, $cellRef) {
    $cell = $ws.Range($cellRef)
    if ($cell.Value2() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

function Add-HandbackLink($ws, $destCell, $displayText, $targetUrl) {
    $cell = $ws.Range($destCell)
    $cell.Value = $displayText
    $ws.Hyperlinks.Add($cell, $targetUrl, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
    $cell.Font.Underline = 2
    $cell.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# 1. Overview sheet: status column for both languages, both tracked files
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
Set-HandbackStatus $overview "B2"
Set-HandbackStatus $overview "C2"
Set-HandbackStatus $overview "B3"
Set-HandbackStatus $overview "C3"

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackStatus $zhcn "B2"
Set-HandbackStatus $zhcn "B3"

$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a6091c99208253e9987525c67ca0329750387a9e/e2e/d927d043-cb01-4ba3-b9f3-42dc29acfa14.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1fddfb5204fa18e12e4a70a0ab23ff728389d2b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d927d043-cb01-4ba3-b9f3-42dc29acfa14.5d2d66fae024143bd68994d27458dcae606c1f88.zh-cn.xlf"

Add-HandbackLink $zhcn "E2" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.md" $zhMdUrl
Add-HandbackLink $zhcn "F2" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.5d2d66fae024143bd68994d27458dcae606c1f88.zh-cn.xlf" $zhXlfUrl
Add-HandbackLink $zhcn "E3" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.md" $zhMdUrl
Add-HandbackLink $zhcn "F3" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.5d2d66fae024143bd68994d27458dcae606c1f88.zh-cn.xlf" $zhXlfUrl

$zhcn.Range("G2").Value = "2016-03-04 03:53:45"
$zhcn.Range("G3").Value = "2016-03-04 03:53:45"

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
Set-HandbackStatus $dede "B2"
Set-HandbackStatus $dede "B3"

$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a6091c99208253e9987525c67ca0329750387a9e/e2e/d927d043-cb01-4ba3-b9f3-42dc29acfa14.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/31b655d497ed5ccb011b57a8226733f79bc5de14/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/d927d043-cb01-4ba3-b9f3-42dc29acfa14.5d2d66fae024143bd68994d27458dcae606c1f88.de-de.xlf"

Add-HandbackLink $dede "E2" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.md" $deMdUrl
Add-HandbackLink $dede "F2" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.5d2d66fae024143bd68994d27458dcae606c1f88.de-de.xlf" $deXlfUrl
Add-HandbackLink $dede "E3" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.md" $deMdUrl
Add-HandbackLink $dede "F3" "d927d043-cb01-4ba3-b9f3-42dc29acfa14.5d2d66fae024143bd68994d27458dcae606c1f88.de-de.xlf" $deXlfUrl

$dede.Range("G2").Value = "2016-03-04 03:54:12"
$dede.Range("G3").Value = "2016-03-04 03:54:12"
